# "Changes without switch case"
#
# The sheet had two sparse rows (A9 and A15) both holding the shared string
# "mytestdata1", separated by blank rows. The edit:
#   - collapses those blank rows so the two data rows become contiguous
#     rows 5 and 6 right after the existing rows 1-4,
#   - replaces their text (row 5 keeps "https://google.com", row 6 becomes
#     "https://www.bing.com/"),
#   - adds a second column value next to each ("Text3" / "Text4"),
#   - moves the active selection to A9,
#   - widens columns A and B to fit the new, longer content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull rows 9 and 15 up into rows 5 and 6 by removing the empty rows
# between them (this also shrinks the used range down to A1:B6).
$ws.Range("A5:A8").EntireRow.Delete()
$ws.Range("A6:A10").EntireRow.Delete()

# Write the new cell contents. Column A keeps/changes the url, column B
# gains a brand-new "TextN" label alongside it.
$ws.Range("A6").Value = "https://www.bing.com/"
$ws.Range("B5").Value = "Text3"
$ws.Range("B6").Value = "Text4"
$ws.Range("A5").Value = "https://google.com"

# Resize the columns to fit the new (longer) values.
$ws.Columns.Item(1).ColumnWidth = 16.0533854166667
$ws.Columns.Item(2).ColumnWidth = 16.7213541666667

# Move the selection, matching the post-edit view state.
[void]$ws.Range("A9").Select()

Write-Output "done"
